# Dev Guide update:
#  1. "AppendiX" (TOC list entry)         -> "AppendiCES"
#  2. "8. APPENDIX" (Heading1)            -> "8. APPENDICES"
#  3. Drop the redundant word " guide" from "This guide will also ..."
#     and relocate the "_GoBack" bookmark to the spot that edit left behind
#     (it used to sit in front of the "2.3 Application Files" heading tab).
#  4. Footer page-number field's cached result "5" -> "3"

$d = $word.ActiveDocument

# 1. Known Issues and Future Work -> Appendices, TOC-style list entry.
$d.Content.Find.Execute("AppendiX", $true, $false, $false, $false, $false, $true, 1, $false, "AppendiCES", 2)

# 2. Section 8 heading.
$d.Content.Find.Execute("8. APPENDIX", $true, $false, $false, $false, $false, $true, 1, $false, "8. APPENDICES", 2)

# 3. Locate "This guide will also" in the introduction paragraph, strip the
#    " guide" and drop the _GoBack bookmark exactly where the edit happened
#    (right after "This"). Adding a bookmark with an existing name moves it,
#    so the stale copy up by "2.3" disappears automatically.
$r = $d.Content
$r.Find.Execute("This guide will also", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$editPoint = $r.Start + 4   # length of "This"

$bmSpot = $d.Range($editPoint, $editPoint)
$d.Bookmarks.Add("_GoBack", $bmSpot)

$guideWord = $d.Range($editPoint, $editPoint + 6)   # " guide"
$guideWord.Delete()

# 4. Footer "PAGE" field's cached display text.
$footer = $d.Sections.First.Footers.Item(1)
$pageNum = $footer.Range.Duplicate
$pageNum.Start = 0
$pageNum.End = 1
$pageNum.Text = "3"
